# Applies row-content swaps identified from the diff: rows are re-paired
# (7<->8, 9<->10, 15/16/17 rotate, 21/22/23 rotate, 25<->26, 33<->34).
# Only the row-number-bearing position stays fixed; every data column moves
# with the record. Read-all-then-write-all avoids clobbering sources that
# are also destinations (e.g. the 3-way rotations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","M","P","Q","R","S","T","U","V","W","AC","AD","AE","AG","AW","AX")

# --- snapshot current ("before") values for every row that participates ---
$rows = @(7,8,9,10,15,16,17,21,22,23,25,26,33,34)
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# --- destination row -> source row (content provider) ---
$mapping = @{
    7 = 8
    8 = 7
    9 = 10
    10 = 9
    15 = 17
    16 = 15
    17 = 16
    21 = 23
    22 = 21
    23 = 22
    25 = 26
    26 = 25
    33 = 34
    34 = 33
}

# --- write the swapped content back out ---
foreach ($dst in $mapping.Keys) {
    $src = $mapping[$dst]
    $rowData = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$dst").Value = $rowData[$c]
    }
}

Write-Host "Row swap complete."
